{"js": "// Map of old division expressions to new ones, as described by the diff.\nconst replacements = {\n  \"496\u00f74=\": \"692\u00f76=\",\n  \"677\u00f74=\": \"713\u00f72=\",\n  \"740\u00f76=\": \"831\u00f74=\",\n  \"288\u00f77=\": \"447\u00f79=\",\n  \"728\u00f76=\": \"312\u00f77=\",\n  \"365\u00f77=\": \"273\u00f72=\",\n  \"534\u00f75=\": \"619\u00f78=\",\n  \"928\u00f73=\": \"780\u00f72=\",\n  \"733\u00f73=\": \"528\u00f78=\",\n  \"715\u00f74=\": \"954\u00f79=\",\n  \"638\u00f78=\": \"987\u00f77=\",\n  \"797\u00f74=\": \"365\u00f78=\",\n  \"202\u00f75=\": \"785\u00f76=\",\n  \"132\u00f72=\": \"265\u00f77=\",\n  \"474\u00f77=\": \"904\u00f78=\",\n  \"809\u00f74=\": \"341\u00f72=\",\n  \"127\u00f74=\": \"541\u00f78=\",\n  \"701\u00f77=\": \"249\u00f77=\",\n  \"941\u00f74=\": \"319\u00f77=\",\n  \"880\u00f72=\": \"176\u00f74=\",\n  \"554\u00f73=\": \"108\u00f77=\",\n  \"384\u00f76=\": \"732\u00f79=\",\n  \"555\u00f78=\": \"903\u00f75=\",\n  \"145\u00f73=\": \"159\u00f74=\",\n  \"554\u00f72=\": \"759\u00f72=\",\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, text)) {\n    paragraph.getRange().insertText(replacements[text], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old division expressions to new ones, as described by the diff.\n$map = [ordered]@{\n    \"496\u00f74=\" = \"692\u00f76=\"\n    \"677\u00f74=\" = \"713\u00f72=\"\n    \"740\u00f76=\" = \"831\u00f74=\"\n    \"288\u00f77=\" = \"447\u00f79=\"\n    \"728\u00f76=\" = \"312\u00f77=\"\n    \"365\u00f77=\" = \"273\u00f72=\"\n    \"534\u00f75=\" = \"619\u00f78=\"\n    \"928\u00f73=\" = \"780\u00f72=\"\n    \"733\u00f73=\" = \"528\u00f78=\"\n    \"715\u00f74=\" = \"954\u00f79=\"\n    \"638\u00f78=\" = \"987\u00f77=\"\n    \"797\u00f74=\" = \"365\u00f78=\"\n    \"202\u00f75=\" = \"785\u00f76=\"\n    \"132\u00f72=\" = \"265\u00f77=\"\n    \"474\u00f77=\" = \"904\u00f78=\"\n    \"809\u00f74=\" = \"341\u00f72=\"\n    \"127\u00f74=\" = \"541\u00f78=\"\n    \"701\u00f77=\" = \"249\u00f77=\"\n    \"941\u00f74=\" = \"319\u00f77=\"\n    \"880\u00f72=\" = \"176\u00f74=\"\n    \"554\u00f73=\" = \"108\u00f77=\"\n    \"384\u00f76=\" = \"732\u00f79=\"\n    \"555\u00f78=\" = \"903\u00f75=\"\n    \"145\u00f73=\" = \"159\u00f74=\"\n    \"554\u00f72=\" = \"759\u00f72=\"\n}\n\nforeach ($key in $map.Keys) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($key, $false, $false, $false, $false, $false, $true, 1, $false, $map[$key], 2)\n}\n"}
